$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","E","F","G","H","I","J","L","O")

$data = @(
    @(16.95397778511553,12.41591857464694,16.83251845988741,33.78795527274045,23.86155347129156,13.16584197284808,20.71089771946079,7.089095426311948,13.27791223070825,19.3678682405579),
    @(16.23233744172688,12.2744208174767,16.87193062414744,33.87476005805679,24.01687568126194,13.22897922036622,20.86122561894723,7.077058354459258,13.22256874017467,19.48091336500062),
    @(15.77196195832255,12.18724561353893,16.89871724618124,33.93741293703062,24.12287635178147,13.27025370527841,20.95846025189482,7.069785763098004,13.18998107027397,19.55552147312349),
    @(15.58023193902234,12.15167670948512,16.91028366992447,33.96528975616403,24.16872572203554,13.28770415861398,20.99932651930513,7.066852661713694,13.17706022786256,19.58722902946585),
    @(15.5481534121687,12.14576871977986,16.91224356422918,33.97006011168806,24.17649862452846,13.29063989376014,21.00618742729424,7.066367490348995,13.17493664945361,19.59257272799871),
    @(15.76939259757715,12.18676605867744,16.89887060050421,33.93777940741028,24.12348397664709,13.27048649407893,20.95900635570298,7.06974608155906,13.18980535115677,19.55594381570582),
    @(16.70886842047007,12.36720925435288,16.84557079390611,33.81594006621059,23.91288973011612,13.18709123773657,20.76170819885991,7.084920996119864,13.25854579001166,19.40576560129402),
    @(18.40569364908223,12.71741655001105,16.76157738229819,33.65151939630942,23.58516915937983,13.04344743012244,20.41385395674999,7.115581800184427,13.40401412238187,19.15264153168335),
    @(19.55434066157878,12.97077847385733,16.71237674773034,33.57648349993975,23.39755918437957,12.95003253008576,20.18196016462431,7.138612187913547,13.51685425511291,18.99207646619586),
    @(20.05409612426459,13.08483077572763,16.69270856388507,33.55234209929023,23.32400088039162,12.91016638890749,20.08157993009388,7.149188591331187,13.56936222508939,18.92458808012002),
    @(20.23996547262964,13.1278175680954,16.6856506959637,33.54464058737363,23.29786053471173,12.89544818809562,20.04430150990249,7.153207029694999,13.58940411801409,18.89983411331324),
    @(20.20008685003074,13.11856905005417,16.68715338952126,33.54623514365369,23.30341376258799,12.89860118961628,20.05229748705082,7.152341007312654,13.58508087489825,18.90512957335341),
    @(20.06945592282675,13.08837154508527,16.69212009273626,33.55167960924445,23.32181582506338,12.90894793172427,20.07849832102259,7.149518928002112,13.57100797507487,18.92253545247585),
    @(19.98899812339438,13.06984749521965,16.69521313229407,33.55520215706602,23.33331148297393,12.91533487555639,20.09464256002462,7.147792037120541,13.56240820084082,18.93330166977722),
    @(19.52121309916786,12.96329815042865,16.71371670010413,33.57826258582216,23.40260526613645,12.95269079246317,20.18862294787662,7.137922920685252,13.51344548155619,18.99659906765149),
    @(19.2283316222604,12.89760427716916,16.72576293084605,33.59497148497073,23.44814778096929,12.97628088009699,20.24758454305335,7.131893409638572,13.48370206536731,19.03685505920897),
    @(19.05773574829873,12.85970734220143,16.73294701704264,33.60552242742173,23.47545088498382,12.99009669078391,20.28197869006274,7.128434843150772,13.46670583797718,19.06053172603239),
    @(18.99961109175785,12.84685783737751,16.73542329603194,33.60925619714376,23.48488497887469,12.9948169744278,20.29370661084402,7.127265483880607,13.46097066878678,19.06863785999202),
    @(19.2597314488406,12.90460926268364,16.72445415560068,33.59309543640154,23.44318485739596,12.97374406678329,20.24125820909374,7.132534290388917,13.48685685170122,19.03251563730025),
    @(20.10791778792143,13.09724700514357,16.69065066725898,33.55004132578331,23.3163640061344,12.90589857541454,20.07078260702396,7.150347485606591,13.57513731963657,18.91740111550793),
    @(20.64253508740581,13.22195452907884,16.6708315181179,33.53029913189468,23.24348181792658,12.86376243463377,19.96364117919402,7.162066993825114,13.63375098709316,18.84684606364884),
    @(20.35903415277404,13.15551434191435,16.68120141714084,33.54006678934262,23.28145884684164,12.88604946220719,20.02043390216455,7.155805296954178,13.60238749654713,18.88407321550937),
    @(19.24554247667439,12.90144270722455,16.72504504753068,33.59394065533188,23.44542510759013,12.97489017028116,20.2441167984847,7.132244523580286,13.48543024909271,19.03447583008137),
    @(17.96328768562276,12.62324490230288,16.78210253333732,33.6879846277152,23.66457430779393,13.08017800495186,20.50379114089024,7.107195321765081,13.36357333519193,19.21667028247492)
)

$startRow = 2
for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    $excelRow = $startRow + $r
    for ($c = 0; $c -lt $cols.Length; $c++) {
        $addr = "$($cols[$c])$excelRow"
        $ws.Range($addr).Value = $rowVals[$c]
    }
}